$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 2.75
$ws.Range("S3").Value = 1.5
$ws.Range("T3").Value = 2.5
$ws.Range("U3").Value = 2.5
$ws.Range("V3").Value = 1.5
$ws.Range("W3").Value = 5
$ws.Range("AC3").Value = 7.5
$ws.Range("AE3").Value = 26
$ws.Range("AL3").Value = 67
$ws.Range("AS3").Value = 251
$ws.Range("AT3").Value = 2.5
$ws.Range("AY3").Value = 51

# Row 4 updates
$ws.Range("I4").Value = 2.45
$ws.Range("J4").Value = 3.25
$ws.Range("L4").Value = 3
$ws.Range("P4").Value = 2.95
$ws.Range("Q4").Value = 1.93
$ws.Range("R4").Value = 1.78
$ws.Range("W4").Value = 9.25
$ws.Range("X4").Value = 15
$ws.Range("Y4").Value = 10
$ws.Range("AA4").Value = 23
$ws.Range("AB4").Value = 30
$ws.Range("AC4").Value = 9.25
$ws.Range("AD4").Value = 6
$ws.Range("AH4").Value = 8.25
$ws.Range("AI4").Value = 12.5
$ws.Range("AL4").Value = 20
$ws.Range("AM4").Value = 29
$ws.Range("AO4").Value = 14.5
$ws.Range("AP4").Value = 19.5
$ws.Range("AQ4").Value = 65
$ws.Range("AR4").Value = 90
$ws.Range("AS4").Value = 200
$ws.Range("AU4").Value = 6.4
$ws.Range("AY4").Value = 19
$ws.Range("BA4").Value = 80

$wb.Save()
